$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shared string text constants used below
$S25 = "B"
$S26 = "Boisson sucre"
$S27 = "C"
$S28 = "Compote sucre"
$S32 = "CS"
$S33 = "Compote sel"
$S598 = "NF EN 17444"
$S817 = "Protéalpes"
$S818 = "Boisson Endurance Fruits Rouges ou Menthe"
$S819 = "Ultimum "
$S820 = "Ultimum Sport Banane"
$S821 = "Ultimum Sport Figue"
$S822 = "Ultimum Sport Pruneau"
$S823 = "Ultimum Sport Abricot"
$S824 = "Ultimum Sport Date"
$S825 = "Ultimum Sport Ananas"
$S826 = "Ultimum Sport Mangue"
$S827 = "Ultimum Sport Kiwi ou Pêche  "
$S828 = "Ultimum Sport Mix Pruneau Cranberry ou Banane Cassis ou Date Goji"
$S829 = "Ultimum Oxygen Datte ou Pruneaux ou Abricot"
$S830 = "Ultimum Oxygen Figue"
$S831 = "Ultimum Oxygen Mangue"
$S832 = "Ultimum Mix Salé Petits Pois Kiwi"
$S833 = "Ultimum Mix Salé Tomate Carotte"

# --- Row 695 ---
$ws.Range("A695").Value = $S817
$ws.Range("B695").Value = $S25
$ws.Range("C695").Value = $S26
$ws.Range("D695").Value = $S818
$ws.Range("E695").Value = 1
$ws.Range("F695").Value = 91.5
$ws.Range("G695").Formula = "=(E695/100)*F695"
$ws.Range("H695").Value = 0
$ws.Range("I695").Formula = "=(E695/100)*H695"
$ws.Range("J695").Value = 0
$ws.Range("K695").Formula = "=G695/E695"
$ws.Range("L695").Formula = "=0.576/100"
$ws.Range("M695").Formula = "=(27.95/800)/G695"
$ws.Range("N695").Value = 1
$ws.Range("O695").Value = 0
$ws.Range("P695").Value = 0
$ws.Range("Q695").Value = 0
$ws.Range("U695").Value = $S598
$ws.Range("V695").Value = 1

# --- Row 696 ---
$ws.Range("A696").Value = $S819
$ws.Range("B696").Value = $S27
$ws.Range("C696").Value = $S28
$ws.Range("D696").Value = $S820
$ws.Range("E696").Value = 70
$ws.Range("F696").Value = 59
$ws.Range("G696").Formula = "=(E696/100)*F696"
$ws.Range("H696").Value = 0.6
$ws.Range("I696").Formula = "=(E696/100)*H696"
$ws.Range("J696").Value = 0
$ws.Range("K696").Formula = "=G696/E696"
$ws.Range("L696").Value = 0
$ws.Range("M696").Formula = "=2.5/G696"
$ws.Range("N696").Value = 0
$ws.Range("O696").Value = 0
$ws.Range("P696").Value = 0
$ws.Range("Q696").Value = 0
$ws.Range("V696").Value = 0

# --- Row 697 ---
$ws.Range("A697").Value = $S819
$ws.Range("B697").Value = $S27
$ws.Range("C697").Value = $S28
$ws.Range("D697").Value = $S821
$ws.Range("E697").Value = 70
$ws.Range("F697").Formula = "=52/0.7"
$ws.Range("G697").Formula = "=(E697/100)*F697"
$ws.Range("H697").Formula = "=1.1/0.7"
$ws.Range("I697").Formula = "=(E697/100)*H697"
$ws.Range("J697").Value = 0
$ws.Range("K697").Formula = "=G697/E697"
$ws.Range("L697").Value = 0
$ws.Range("M697").Formula = "=2.5/G697"
$ws.Range("N697").Value = 0
$ws.Range("O697").Value = 0
$ws.Range("P697").Value = 0
$ws.Range("Q697").Value = 0
$ws.Range("V697").Value = 0

# --- Row 698 ---
$ws.Range("A698").Value = $S819
$ws.Range("B698").Value = $S27
$ws.Range("C698").Value = $S28
$ws.Range("D698").Value = $S822
$ws.Range("E698").Value = 70
$ws.Range("F698").Formula = "=44/0.7"
$ws.Range("G698").Formula = "=(E698/100)*F698"
$ws.Range("H698").Formula = "=0.7/0.7"
$ws.Range("I698").Formula = "=(E698/100)*H698"
$ws.Range("J698").Value = 0
$ws.Range("K698").Formula = "=G698/E698"
$ws.Range("L698").Value = 0
$ws.Range("M698").Formula = "=2.5/G698"
$ws.Range("N698").Value = 0
$ws.Range("O698").Value = 0
$ws.Range("P698").Value = 0
$ws.Range("Q698").Value = 0
$ws.Range("V698").Value = 0

# --- Row 699 ---
$ws.Range("A699").Value = $S819
$ws.Range("B699").Value = $S27
$ws.Range("C699").Value = $S28
$ws.Range("D699").Value = $S823
$ws.Range("E699").Value = 70
$ws.Range("F699").Formula = "=47/0.7"
$ws.Range("G699").Formula = "=(E699/100)*F699"
$ws.Range("H699").Formula = "=1.2/0.7"
$ws.Range("I699").Formula = "=(E699/100)*H699"
$ws.Range("J699").Value = 0
$ws.Range("K699").Formula = "=G699/E699"
$ws.Range("L699").Value = 0
$ws.Range("M699").Formula = "=2.5/G699"
$ws.Range("N699").Value = 0
$ws.Range("O699").Value = 0
$ws.Range("P699").Value = 0
$ws.Range("Q699").Value = 0
$ws.Range("V699").Value = 0

# --- Row 700 ---
$ws.Range("A700").Value = $S819
$ws.Range("B700").Value = $S27
$ws.Range("C700").Value = $S28
$ws.Range("D700").Value = $S824
$ws.Range("E700").Value = 70
$ws.Range("F700").Formula = "=39/0.7"
$ws.Range("G700").Formula = "=(E700/100)*F700"
$ws.Range("H700").Formula = "=1/0.7"
$ws.Range("I700").Formula = "=(E700/100)*H700"
$ws.Range("J700").Value = 0
$ws.Range("K700").Formula = "=G700/E700"
$ws.Range("L700").Value = 0
$ws.Range("M700").Formula = "=2.5/G700"
$ws.Range("N700").Value = 0
$ws.Range("O700").Value = 0
$ws.Range("P700").Value = 0
$ws.Range("Q700").Value = 0
$ws.Range("V700").Value = 0

# --- Row 701 ---
$ws.Range("A701").Value = $S819
$ws.Range("B701").Value = $S27
$ws.Range("C701").Value = $S28
$ws.Range("D701").Value = $S825
$ws.Range("E701").Value = 70
$ws.Range("F701").Formula = "=52/0.7"
$ws.Range("G701").Formula = "=(E701/100)*F701"
$ws.Range("H701").Formula = "=0.4/0.7"
$ws.Range("I701").Formula = "=(E701/100)*H701"
$ws.Range("J701").Value = 0
$ws.Range("K701").Formula = "=G701/E701"
$ws.Range("L701").Formula = "=0.04*0.4"
$ws.Range("M701").Formula = "=2.5/G701"
$ws.Range("N701").Value = 0
$ws.Range("O701").Value = 0
$ws.Range("P701").Value = 0
$ws.Range("Q701").Value = 0
$ws.Range("V701").Value = 0

# --- Row 702 ---
$ws.Range("A702").Value = $S819
$ws.Range("B702").Value = $S27
$ws.Range("C702").Value = $S28
$ws.Range("D702").Value = $S826
$ws.Range("E702").Value = 70
$ws.Range("F702").Formula = "=54.7/0.7"
$ws.Range("G702").Formula = "=(E702/100)*F702"
$ws.Range("H702").Formula = "=0.8/0.7"
$ws.Range("I702").Formula = "=(E702/100)*H702"
$ws.Range("J702").Value = 0
$ws.Range("K702").Formula = "=G702/E702"
$ws.Range("L702").Formula = "=0.35*0.4"
$ws.Range("M702").Formula = "=2.5/G702"
$ws.Range("N702").Value = 0
$ws.Range("O702").Value = 0
$ws.Range("P702").Value = 0
$ws.Range("Q702").Value = 0
$ws.Range("V702").Value = 0

# --- Row 703 ---
$ws.Range("A703").Value = $S819
$ws.Range("B703").Value = $S27
$ws.Range("C703").Value = $S28
$ws.Range("D703").Value = $S827
$ws.Range("E703").Value = 70
$ws.Range("F703").Formula = "=40/0.7"
$ws.Range("G703").Formula = "=(E703/100)*F703"
$ws.Range("H703").Formula = "=0.4/0.7"
$ws.Range("I703").Formula = "=(E703/100)*H703"
$ws.Range("J703").Value = 0
$ws.Range("K703").Formula = "=G703/E703"
$ws.Range("L703").Value = 0.016
$ws.Range("M703").Formula = "=2.5/G703"
$ws.Range("N703").Value = 0
$ws.Range("O703").Value = 0
$ws.Range("P703").Value = 0
$ws.Range("Q703").Value = 0
$ws.Range("V703").Value = 0

# --- Row 704 ---
$ws.Range("A704").Value = $S819
$ws.Range("B704").Value = $S27
$ws.Range("C704").Value = $S28
$ws.Range("D704").Value = $S828
$ws.Range("E704").Value = 70
$ws.Range("F704").Formula = "=40/0.7"
$ws.Range("G704").Formula = "=(E704/100)*F704"
$ws.Range("H704").Formula = "=0.5/0.7"
$ws.Range("I704").Formula = "=(E704/100)*H704"
$ws.Range("J704").Value = 0
$ws.Range("K704").Formula = "=G704/E704"
$ws.Range("L704").Value = 0.01
$ws.Range("M704").Formula = "=2.5/G704"
$ws.Range("N704").Value = 0
$ws.Range("O704").Value = 0
$ws.Range("P704").Value = 0
$ws.Range("Q704").Value = 0
$ws.Range("V704").Value = 0

# --- Row 705 ---
$ws.Range("A705").Value = $S819
$ws.Range("B705").Value = $S27
$ws.Range("C705").Value = $S28
$ws.Range("D705").Value = $S829
$ws.Range("E705").Value = 70
$ws.Range("F705").Formula = "=20.3/0.7"
$ws.Range("G705").Formula = "=(E705/100)*F705"
$ws.Range("H705").Formula = "=0.56/0.7"
$ws.Range("I705").Formula = "=(E705/100)*H705"
$ws.Range("J705").Value = 0
$ws.Range("K705").Formula = "=G705/E705"
$ws.Range("L705").Value = 0
$ws.Range("M705").Formula = "=2.5/G705"
$ws.Range("N705").Value = 0
$ws.Range("O705").Value = 0
$ws.Range("P705").Value = 0
$ws.Range("Q705").Value = 0
$ws.Range("V705").Value = 0

# --- Row 706 ---
$ws.Range("A706").Value = $S819
$ws.Range("B706").Value = $S27
$ws.Range("C706").Value = $S28
$ws.Range("D706").Value = $S830
$ws.Range("E706").Value = 70
$ws.Range("F706").Formula = "=16.59/0.7"
$ws.Range("G706").Formula = "=(E706/100)*F706"
$ws.Range("H706").Formula = "=1.05/0.7"
$ws.Range("I706").Formula = "=(E706/100)*H706"
$ws.Range("J706").Value = 0
$ws.Range("K706").Formula = "=G706/E706"
$ws.Range("L706").Value = 0
$ws.Range("M706").Formula = "=2.5/G706"
$ws.Range("N706").Value = 0
$ws.Range("O706").Value = 0
$ws.Range("P706").Value = 0
$ws.Range("Q706").Value = 0
$ws.Range("V706").Value = 0

# --- Row 707 ---
$ws.Range("A707").Value = $S819
$ws.Range("B707").Value = $S27
$ws.Range("C707").Value = $S28
$ws.Range("D707").Value = $S831
$ws.Range("E707").Value = 70
$ws.Range("F707").Formula = "=25.3/0.7"
$ws.Range("G707").Formula = "=(E707/100)*F707"
$ws.Range("H707").Formula = "=0.8/0.7"
$ws.Range("I707").Formula = "=(E707/100)*H707"
$ws.Range("J707").Value = 0
$ws.Range("K707").Formula = "=G707/E707"
$ws.Range("L707").Value = 0
$ws.Range("M707").Formula = "=2.5/G707"
$ws.Range("N707").Value = 0
$ws.Range("O707").Value = 0
$ws.Range("P707").Value = 0
$ws.Range("Q707").Value = 0
$ws.Range("V707").Value = 0

# --- Row 708 ---
$ws.Range("A708").Value = $S819
$ws.Range("B708").Value = $S32
$ws.Range("C708").Value = $S33
$ws.Range("D708").Value = $S832
$ws.Range("E708").Value = 70
$ws.Range("F708").Formula = "=6.72/0.7"
$ws.Range("G708").Formula = "=(E708/100)*F708"
$ws.Range("H708").Formula = "=2.52/0.7"
$ws.Range("I708").Formula = "=(E708/100)*H708"
$ws.Range("J708").Value = 0
$ws.Range("K708").Formula = "=G708/E708"
$ws.Range("L708").Formula = "=0.42*0.4"
$ws.Range("M708").Formula = "=2.5/G708"
$ws.Range("N708").Value = 0
$ws.Range("O708").Value = 0
$ws.Range("P708").Value = 0
$ws.Range("Q708").Value = 0
$ws.Range("V708").Value = 0

# --- Row 709 ---
$ws.Range("A709").Value = $S819
$ws.Range("B709").Value = $S32
$ws.Range("C709").Value = $S33
$ws.Range("D709").Value = $S833
$ws.Range("E709").Value = 70
$ws.Range("F709").Formula = "=3.01/0.7"
$ws.Range("G709").Formula = "=(E709/100)*F709"
$ws.Range("H709").Formula = "=0.63/0.7"
$ws.Range("I709").Formula = "=(E709/100)*H709"
$ws.Range("J709").Value = 0
$ws.Range("K709").Formula = "=G709/E709"
$ws.Range("L709").Formula = "=0.49*0.4"
$ws.Range("M709").Formula = "=2.5/G709"
$ws.Range("N709").Value = 0
$ws.Range("O709").Value = 0
$ws.Range("P709").Value = 0
$ws.Range("Q709").Value = 0
$ws.Range("V709").Value = 0

# --- Blank placeholder rows 710-714 (mirrors pre-existing blank K-only rows) ---
$ws.Range("K710").Formula = "=G710/E710"
$ws.Range("K711").Formula = "=G711/E711"
$ws.Range("K712").Formula = "=G712/E712"
$ws.Range("K713").Formula = "=G713/E713"
$ws.Range("K714").Formula = "=G714/E714"

# --- Reset D685:D687 to default (general) style, matching removal of redundant font/style ---
$ws.Range("D685").NumberFormat = "General"
$ws.Range("D686").NumberFormat = "General"
$ws.Range("D687").NumberFormat = "General"

# --- Final selection state ---
$ws.Range("X707").Select()
